$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current row 38 (pushes old rows 38-61 down to 41-64)
$ws.Range("A38:A40").EntireRow.Insert()

# Common/static values shared by every data row in this sheet
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad  = "Sin especificar"
$unidad    = "$/bandeja 3 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnidad  = 3

function Set-DataRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-DataRow 38 45079 "Primera" 200 9000 10000 9500 3167
Set-DataRow 39 45079 "Segunda" 100 7000 8000 7500 2500
Set-DataRow 40 45079 "Tercera" 160 4000 5000 4500 1500
